# api for spam classification using count vectorizer
# - Rename the "Sales Category" sheet to "Config_Sheet"
# - Bump the version value in B1 from 3 to 6
# - Move the active selection to B1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales Category")

$ws.Name = "Config_Sheet"

$ws.Range("B1").Value = 6

$ws.Range("B1").Select() | Out-Null
